$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-of dates) ---
$ws.Range("A8").Characters(21, 2).Text = "45"
$ws.Range("C9").Characters(48, 9).Text = "11/10/2024"
$ws.Range("C9").Characters(27, 10).Text = "11/4/2024"

# --- Style changes: copy number-format from a stable donor cell, then set new value ---
# donors: D14 = "0"-text count style (13), E14 = "***.*"-text pct style (13),
#         J14 = plain count style (14), K15 = plain pct style (15)
$ws.Range("J14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1
$ws.Range("J14").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1
$ws.Range("D14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("G27").Value = "0"
$ws.Range("E14").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("H27").Value = "***.*"
$ws.Range("D14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = "0"
$ws.Range("J14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1
$ws.Range("K15").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("D14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = "0"
$ws.Range("J14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1
$ws.Range("K15").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$excel.CutCopyMode = 0

# --- Plain value updates (no style change) ---
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 4
$ws.Range("I14").Value = 8
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 60
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = -79.487179487179
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -71.428571428571
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = -44.444444444444
$ws.Range("I16").Value = 244
$ws.Range("J16").Value = 290
$ws.Range("K16").Value = -15.862068965517
$ws.Range("L16").Value = -7.575757575757
$ws.Range("M16").Value = -0.813008130081
$ws.Range("N16").Value = -67.248322147651
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 28.571428571428
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 38
$ws.Range("H17").Value = -13.157894736842
$ws.Range("I17").Value = 386
$ws.Range("J17").Value = 419
$ws.Range("K17").Value = -7.875894988066
$ws.Range("L17").Value = 2.659574468085
$ws.Range("M17").Value = 14.540059347181
$ws.Range("N17").Value = -14.031180400890
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -12.5
$ws.Range("I18").Value = 164
$ws.Range("J18").Value = 155
$ws.Range("K18").Value = 5.806451612903
$ws.Range("L18").Value = 19.708029197080
$ws.Range("M18").Value = -13.684210526315
$ws.Range("N18").Value = -82.231852654387
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -62.5
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 373
$ws.Range("J19").Value = 394
$ws.Range("K19").Value = -5.329949238578
$ws.Range("L19").Value = 5.070422535211
$ws.Range("M19").Value = 60.775862068965
$ws.Range("N19").Value = -4.113110539845
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 6.25
$ws.Range("I20").Value = 155
$ws.Range("J20").Value = 264
$ws.Range("K20").Value = -41.287878787878
$ws.Range("L20").Value = -14.835164835164
$ws.Range("M20").Value = 124.63768115942
$ws.Range("N20").Value = -57.880434782608
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -34.375
$ws.Range("F21").Value = 113
$ws.Range("G21").Value = 127
$ws.Range("H21").Value = -11.023622047244
$ws.Range("I21").Value = 1364
$ws.Range("J21").Value = 1550
$ws.Range("K21").Value = -12
$ws.Range("L21").Value = 1.262063845582
$ws.Range("M21").Value = 24.225865209471
$ws.Range("N21").Value = -53.621217273036
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 18
$ws.Range("K22").Value = -18.181818181818
$ws.Range("L22").Value = -30.769230769230
$ws.Range("M22").Value = -18.181818181818
$ws.Range("L23").Value = 43.75
$ws.Range("M23").Value = 109.090909090909
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 53.333333333333
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = -15.492957746478
$ws.Range("I24").Value = 701
$ws.Range("J24").Value = 792
$ws.Range("K24").Value = -11.489898989899
$ws.Range("L24").Value = 2.936857562408
$ws.Range("M24").Value = 25.178571428571
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -62.5
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = -50
$ws.Range("I25").Value = 191
$ws.Range("J25").Value = 246
$ws.Range("K25").Value = -22.357723577235
$ws.Range("L25").Value = -27.924528301886
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 55
$ws.Range("G26").Value = 52
$ws.Range("H26").Value = 5.769230769230
$ws.Range("I26").Value = 576
$ws.Range("J26").Value = 495
$ws.Range("K26").Value = 16.363636363636
$ws.Range("L26").Value = 22.814498933901
$ws.Range("M26").Value = -0.173310225303
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 300
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = -22.222222222222
$ws.Range("I28").Value = 114
$ws.Range("J28").Value = 89
$ws.Range("K28").Value = 28.089887640449
$ws.Range("L28").Value = 37.349397590361
$ws.Range("J29").Value = 21
$ws.Range("K29").Value = -23.809523809523
$ws.Range("J30").Value = 19
$ws.Range("K30").Value = -26.315789473684
